# Updated run for publication
# Apply refreshed frequency-table values (columns B:X, rows 2:5) onto the
# active worksheet of the already-open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.227544910179641,0.0172155688622754,0.789670658682635,0.0853293413173653,0.697604790419162,0.0868263473053892,0.0658682634730539,0.0232035928143713,0.921407185628742,0.0434131736526946,0.0988023952095808,0.0029940119760479,0.000748502994011976,0.000748502994011976,0.000748502994011976,0,0.852544910179641,0,0.997754491017964,0.907934131736527,0.0172155688622754,0.0254491017964072,0.0359281437125748),
    @(0.00673652694610778,0.157185628742515,0.155688622754491,0.703592814371258,0.114520958083832,0.110778443113772,0.710329341317365,0.0973053892215569,0.0134730538922156,0.00374251497005988,0.678892215568862,0.0134730538922156,0.000748502994011976,0,0.987275449101796,0,0,0.995508982035928,0.00149700598802395,0.0853293413173653,0.0860778443113772,0.226047904191617,0.00374251497005988),
    @(0.0164670658682635,0.820359281437126,0.0291916167664671,0.00598802395209581,0.0187125748502994,0.0134730538922156,0.168413173652695,0.854041916167665,0.0179640718562874,0.780688622754491,0.170658682634731,0.197604790419162,0.00374251497005988,0.462574850299401,0,0.000748502994011976,0.00598802395209581,0.00374251497005988,0.000748502994011976,0.000748502994011976,0.893712574850299,0.709580838323353,0.714820359281437),
    @(0.749251497005988,0.00523952095808383,0.0254491017964072,0.205089820359281,0.169161676646707,0.788922155688623,0.0553892215568862,0.0254491017964072,0.0471556886227545,0.172155688622754,0.0516467065868263,0.785928143712575,0.994760479041916,0.536676646706587,0.0119760479041916,0.999251497005988,0.141467065868263,0.000748502994011976,0,0.00598802395209581,0.00224550898203593,0.0389221556886228,0.245508982035928)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowIndex = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        # Column B is index 2
        $colIndex = $j + 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
